$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated case count for 2021-03-24 (row 393) cascades the running total
# (column B) forward through every later row via the existing shared
# formula, so only the raw input needs to change here.
$ws.Range("C393").Value = 94

# New data entry for 2021-05-17 and 2021-05-18 (rows 447-448).
$ws.Range("C447").Value = 56
$ws.Range("C448").Value = 39

# First-time data entry for 2021-05-19 (row 449), previously blank.
$ws.Range("C449").Value = 2
$ws.Range("L449").Value = "0"
$ws.Range("M449").Value = "0"
